$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the test name for the existing "tc1" row (row 2).
$ws.Cells.Item(2, 2).Value = "testLoginJavascript"

# Update row 3 (tc2): rename the test and reuse the SwatiChetty credentials.
$ws.Cells.Item(3, 2).Value = "testLoginJasmine"
$ws.Cells.Item(3, 3).Value = '{"username":"SwatiChetty","password":"123456"}'

# Add a new row 4 for another tc2 variant (Selenium), reusing the original
# Swati credentials and the Order Create Successfully validation text.
$ws.Cells.Item(4, 1).Value = "tc2"
$ws.Cells.Item(4, 2).Value = "testLoginSelenium"
$ws.Cells.Item(4, 3).Value = '{"username":"Swati","password":"123"}'
$ws.Cells.Item(4, 4).Value = '{"textToValidate":"Order Create Successfully"}'

# Move the active selection to C3, matching the saved workbook state.
$ws.Range("C3").Select()
